$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the field name value in A41
$ws.Range("A41").Value = "qualifiedProfessionalOrganizati"

# Add a comment to A41
$cmt = $ws.Range("A41").AddComment("Roy Jeong:`nNever rename fields. The field name is not incorrect, it is currently actually referenced as such in AGOL internally.")
$cmt.Author = "Roy Jeong"

# Update the view: scroll to show row 26, select B46
$ws.Application.ActiveWindow.ScrollRow = 26
$ws.Range("B46").Select()
